# Update "想去人数" (column F) values on the "展览" (sheet index 1)
# and "全部类型" (sheet index 4) worksheets to match the published
# gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1 / rId1) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 6985
$ws1.Range("F3").Value  = 0
$ws1.Range("F4").Value  = 62
$ws1.Range("F5").Value  = 0
$ws1.Range("F7").Value  = 6875
$ws1.Range("F10").Value = 0
$ws1.Range("F13").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 0
$ws1.Range("F17").Value = 49
$ws1.Range("F19").Value = 0
$ws1.Range("F20").Value = 5245
$ws1.Range("F21").Value = 0
$ws1.Range("F23").Value = 662
$ws1.Range("F25").Value = 0

# --- Sheet "全部类型" (index 4 / rId4) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 6985
$ws4.Range("F3").Value  = 0
$ws4.Range("F4").Value  = 62
$ws4.Range("F5").Value  = 456
$ws4.Range("F7").Value  = 0
$ws4.Range("F9").Value  = 0
$ws4.Range("F12").Value = 109
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 0
$ws4.Range("F16").Value = 416
$ws4.Range("F17").Value = 49
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 5245
$ws4.Range("F22").Value = 0
$ws4.Range("F25").Value = 0
